$d = $word.ActiveDocument

# --- 1. "HTML Code" character style: sz 20 -> 18 (10pt -> 9pt) -------------
$htmlCode = $d.Styles.Item("HTML Code")
$htmlCode.Font.Size = 9

# --- 2. "Code" paragraph style: sz 20 -> 18 (10pt -> 9pt) -------------------
$code = $d.Styles.Item("Code")
$code.Font.Size = 9

# --- 3. "Source Code" paragraph style: add sz 20 (10pt) ---------------------
$sourceCode = $d.Styles.Item("Source Code")
$sourceCode.Font.Size = 10

# --- 4. "Verbatim" paragraph style: add sz 18 (9pt) -------------------------
$verbatim = $d.Styles.Item("Verbatim")
$verbatim.Font.Size = 9

# --- 5. "Code Block" paragraph style: no rPr change (rsid-only in source) --
# (nothing further to do for this style's content)

# --- 6/7. New built-in "annotation text" / "Comment Text Char" styles ------
$commentText = $d.Styles.Add("Comment Text", 1)
$commentText.BaseStyle = "Normal"
$commentText.Priority = 99
$commentText.UnhideWhenUsed = $true
$commentText.Font.Size = 10
$commentText.Font.SizeBi = 10

$commentTextChar = $d.Styles.Add("Comment Text Char", 2)
$commentTextChar.BaseStyle = "Default Paragraph Font"
$commentTextChar.Priority = 99
$commentTextChar.Font.Size = 10
$commentTextChar.Font.SizeBi = 10

$commentText.LinkStyle = $commentTextChar
$commentTextChar.LinkStyle = $commentText
